$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted right above the existing row 67,
# pushing the former rows 67-180 down to 68-181 (the sheet's used range
# grows from A1:R180 to A1:R181). Insert a whole row so everything below
# shifts down automatically, then populate the new row with its data.
$ws.Rows.Item(67).Insert()

$ws.Cells.Item(67, 1).Value  = 3
$ws.Cells.Item(67, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(67, 3).Value  = "Coquimbo"
$ws.Cells.Item(67, 4).Value  = 44580
$ws.Cells.Item(67, 5).Value  = 5
$ws.Cells.Item(67, 6).Value  = 100112010
$ws.Cells.Item(67, 7).Value  = "Achicoria"
$ws.Cells.Item(67, 8).Value  = "Sin especificar"
$ws.Cells.Item(67, 9).Value  = "Primera"
$ws.Cells.Item(67, 10).Value = 50
$ws.Cells.Item(67, 11).Value = 6000
$ws.Cells.Item(67, 12).Value = 6000
$ws.Cells.Item(67, 13).Value = 6000
$ws.Cells.Item(67, 14).Value = '$/caja 16 unidades'
$ws.Cells.Item(67, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(67, 16).Value = 375
$ws.Cells.Item(67, 17).Value = 16
$ws.Cells.Item(67, 18).Value = "Hortaliza"
